# Adapt column header formatting to respective input file names.
# Renames header cells from *_old / *_new suffixes to *_FV2210 / *_FV2304,
# wraps the data range in an Excel Table (ListObject), and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells in row 1 (A1:J1 use "_old" -> "_FV2210", L1:U1 use "_new" -> "_FV2304").
# Column K ("diff") stays unchanged.
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $value = $cell.Value2
    if ($value -like "*_old") {
        $cell.Value = ($value -replace "_old$", "_FV2210")
    } elseif ($value -like "*_new") {
        $cell.Value = ($value -replace "_new$", "_FV2304")
    }
}

# Freeze the header row (row 1).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Wrap the data range (A1:U75) in an Excel Table (ListObject) with headers, matching Table1.
$tableRange = $ws.Range("A1:U75")
$listObject = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = $null
